# Update dSF (column F) values to reflect repulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = -2
$ws.Range("F6").Value = -4
$ws.Range("F8").Value = 1
$ws.Range("F14").Value = 5
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = -6
